# Append four new daily rows (2026-01-17 and 2026-01-18, one per station)
# to the bottom of the existing data table on Sheet1, then move the
# selection down to match where Excel would leave the cursor after
# typing the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows are appended right after the current last row (row 33).
$startRow = 34

$newRows = @(
    @(46039, "四方坪站", 12743.6,            9327.49, 3541.18,            539),
    @(46039, "高岭站",   4460.8900000000003, 3527.9,  1236.8399999999999, 152),
    @(46040, "四方坪站", 11465.654,          9383.41, 3684.83,            491),
    @(46040, "高岭站",   3847.77,            3211.52, 1043.33,            139)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}

# Match the author's final cursor position / scrolled view (H36) after
# entering the new data.
$ws.Range("H36").Select() | Out-Null
